$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: fill in the completion date (C17)
$ws.Range("C17").Value = "2023-01-30"

# Row 18: fill in the completion date (C18)
$ws.Range("C18").Value = "2023-01-31"

# Row 19: new problem entry - "leetcode 202" / "快乐数", completed 2023-02-01
$ws.Range("A19").Value = "leetcode 202"
$ws.Range("B19").Value = "快乐数"
$ws.Range("C19").Value = "2023-02-01"

# Row 20: problem number label, completed 2023-02-01
$ws.Range("A20").Value = "leetcode 202"
$ws.Range("C20").Value = "2023-02-01"

# Row 21: problem number label, completed 2023-02-01
$ws.Range("A21").Value = "leetcode 202"
$ws.Range("C21").Value = "2023-02-01"

# Row 22: problem number label "leetcode"; completion date typed as plain text
$ws.Range("A22").Value = "leetcode"
$ws.Range("C22").Value = "2023-de-01"

# Update the active selection to C22, matching the saved cursor position
$ws.Range("C22").Select()
